# Natmi following Dr Hou advice
# Regenerate the LR-pair table (Cxcl12 -> Cxcr3) with the full cross-product
# of Sending/Target clusters, following Dr Hou's advice.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Cxcl12", "Cxcr3", "FAPs", 3, 1, 198.977211, 596.9316329999999, 0.440791350614085, 0.4407913506140851, 1, 0.3333333333333333, 0.05194133333333333, 0.155824, 0.03788844568234288, 0.03788844568234288, 10.335141642288, 93.01627478059198, 0.01670089914498832, 0.01670089914498832),
    @("ECs", "Cxcl12", "Cxcr3", "M2", 3, 1, 198.977211, 596.9316329999999, 0.440791350614085, 0.4407913506140851, 3, 1, 1.318960333333333, 3.956881, 0.962111554317657, 0.9621115543176572, 262.443048546297, 2361.987436916673, 0.4240904514690967, 0.4240904514690968),
    @("FAPs", "Cxcl12", "Cxcr3", "FAPs", 3, 1, 196.7746736666667, 590.324021, 0.4359121013721307, 0.4359121013721308, 1, 0.3333333333333333, 0.05194133333333333, 0.155824, 0.03788844568234288, 0.03788844568234288, 10.22073891647822, 91.986650248304, 0.01651603197511392, 0.01651603197511392),
    @("FAPs", "Cxcl12", "Cxcr3", "M2", 3, 1, 196.7746736666667, 590.324021, 0.4359121013721307, 0.4359121013721308, 3, 1, 1.318960333333333, 3.956881, 0.962111554317657, 0.9621115543176572, 259.5379891709446, 2335.841902538501, 0.4193960693970167, 0.4193960693970168),
    @("M2", "Cxcl12", "Cxcr3", "FAPs", 2, 0.6666666666666666, 0.158424, 0.475272, 0.0003509544061791369, 0.0003509544061791369, 1, 0.3333333333333333, 0.05194133333333333, 0.155824, 0.03788844568234288, 0.03788844568234288, 0.008228753792, 0.074058784128, 0.00001329711695549713, 0.00001329711695549713),
    @("M2", "Cxcl12", "Cxcr3", "M2", 2, 0.6666666666666666, 0.158424, 0.475272, 0.0003509544061791369, 0.0003509544061791369, 3, 1, 1.318960333333333, 3.956881, 0.962111554317657, 0.9621115543176572, 0.208954971848, 1.880594746632, 0.0003376572892236397, 0.0003376572892236398),
    @("sCs", "Cxcl12", "Cxcr3", "FAPs", 3, 1, 55.49875533333334, 166.496266, 0.122945593607605, 0.122945593607605, 1, 0.3333333333333333, 0.05194133333333333, 0.155824, 0.03788844568234288, 0.03788844568234288, 2.882679350353778, 25.944114153184, 0.004658217445285143, 0.004658217445285143),
    @("sCs", "Cxcl12", "Cxcr3", "M2", 3, 1, 55.49875533333334, 166.496266, 0.122945593607605, 0.122945593607605, 3, 1, 1.318960333333333, 3.956881, 0.962111554317657, 0.9621115543176572, 73.20065683403845, 658.8059115063461, 0.1182873761623198, 0.1182873761623198)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value2 = $rowValues[$j]
    }
}
